# Update "想去人数" (F column) values on the "展览" sheet and the aggregated
# "全部类型" sheet to reflect newly scraped counts.

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Worksheets index 1 / sheet1.xml) ---
$wsExpo = $wb.Worksheets.Item("展览")
$expoUpdates = @{
    2  = 1180
    3  = 1114
    4  = 1923
    6  = 1253
    7  = 68
    8  = 39
    11 = 116
    13 = 817
    14 = 243
    15 = 125
    19 = 226
    20 = 696
    23 = 191
    25 = 906
    27 = 190
}
foreach ($row in $expoUpdates.Keys) {
    $wsExpo.Range("F$row").Value = $expoUpdates[$row]
}

# --- Sheet "全部类型" (aggregated view, sheet4.xml) ---
$wsAll = $wb.Worksheets.Item("全部类型")
$allUpdates = @{
    3  = 1180
    4  = 1114
    5  = 1923
    7  = 1253
    8  = 68
    10 = 39
    13 = 116
    15 = 817
    16 = 243
    17 = 125
    27 = 226
    28 = 696
    31 = 191
    33 = 906
    37 = 190
}
foreach ($row in $allUpdates.Keys) {
    $wsAll.Range("F$row").Value = $allUpdates[$row]
}
